# Generate Report for Handback
# This script updates the zh-cn and de-de localization-status sheets to
# reflect that the handoff packages have now been handed back:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is stamped with the actual handback time
#  - Latest Target File / Latest Handback File columns (F/G) are populated
#    with hyperlinked file names (mirroring the handoff file + the xlf
#    that came back)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet - mirrors the same "Ready for handoff" status text,
# which is shared with the per-language sheets below.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) for both data rows
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Row 2 - 12ca35d4-3a53-4132-b65d-5a124a67de2a
$wsZh.Range("F2").Value = "12ca35d4-3a53-4132-b65d-5a124a67de2a.md"
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", "12ca35d4-3a53-4132-b65d-5a124a67de2a.md")

$wsZh.Range("G2").Value = "12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.zh-cn.xlf"
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/913a889750a8e4e0fbcd5b38739a5194f543e3b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.zh-cn.xlf", "", "", "12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.zh-cn.xlf")

$wsZh.Range("H2").Value = "2016-03-17 18:50:54"

# Row 3 - e05a502e-0b09-4336-88b9-fa363d28dee0
$wsZh.Range("F3").Value = "e05a502e-0b09-4336-88b9-fa363d28dee0.md"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/e05a502e-0b09-4336-88b9-fa363d28dee0.md", "", "", "e05a502e-0b09-4336-88b9-fa363d28dee0.md")

$wsZh.Range("G3").Value = "e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.zh-cn.xlf"
$wsZh.Range("G3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/913a889750a8e4e0fbcd5b38739a5194f543e3b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.zh-cn.xlf", "", "", "e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.zh-cn.xlf")

$wsZh.Range("H3").Value = "2016-03-17 18:50:54"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (C) for both data rows
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Row 2 - 12ca35d4-3a53-4132-b65d-5a124a67de2a
$wsDe.Range("F2").Value = "12ca35d4-3a53-4132-b65d-5a124a67de2a.md"
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/12ca35d4-3a53-4132-b65d-5a124a67de2a.md", "", "", "12ca35d4-3a53-4132-b65d-5a124a67de2a.md")

$wsDe.Range("G2").Value = "12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.de-de.xlf"
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e85420b1dc98b61a4854be8f50dab573519a57b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.de-de.xlf", "", "", "12ca35d4-3a53-4132-b65d-5a124a67de2a.ed3b5d5dd956fa1588ae11ba18bc514cc6a12fd9.de-de.xlf")

$wsDe.Range("H2").Value = "2016-03-17 18:51:00"

# Row 3 - e05a502e-0b09-4336-88b9-fa363d28dee0
$wsDe.Range("F3").Value = "e05a502e-0b09-4336-88b9-fa363d28dee0.md"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/0267841bce86c52a948e038c97a1d398cef05c94/e2e/e05a502e-0b09-4336-88b9-fa363d28dee0.md", "", "", "e05a502e-0b09-4336-88b9-fa363d28dee0.md")

$wsDe.Range("G3").Value = "e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.de-de.xlf"
$wsDe.Range("G3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e85420b1dc98b61a4854be8f50dab573519a57b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.de-de.xlf", "", "", "e05a502e-0b09-4336-88b9-fa363d28dee0.a63e1af62b09e8f69e9ff630316652ca09dbacf0.de-de.xlf")

$wsDe.Range("H3").Value = "2016-03-17 18:51:00"
